$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.823.78"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.271.25"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB (numeric-looking price -> keep as text like source)
$ws.Range("D5").Value = "'303.77"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'92.74"
$ws.Range("E6").Value = "  +0.46%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.530"
$ws.Range("E7").Value = "  +1.83%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'32.63"
$ws.Range("E10").Value = "  +1.48%  "

# Row 11 - OKB
$ws.Range("D11").Value = "'53.61"
$ws.Range("E11").Value = "  -1.10%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "'0.0796"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.42%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.69"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.624.10"
$ws.Range("E15").Value = "  +0.75%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'14.28"
$ws.Range("E16").Value = "  +0.65%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.279.04"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18 - Polygon
$ws.Range("D18").Value = "'0.776"
$ws.Range("E18").Value = "  +2.90%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "41.763.95"
$ws.Range("E19").Value = "  +1.25%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.80"
$ws.Range("E20").Value = "  +4.06%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.40%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'67.17"
$ws.Range("E23").Value = "  +0.32%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "'243.66"
$ws.Range("E24").Value = "  +1.53%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +0.23%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "'1.94"

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.04%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'24.02"
$ws.Range("E28").Value = "  +1.34%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'9.52"
$ws.Range("E29").Value = "  -1.37%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -5.09%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "'35.22"
$ws.Range("E31").Value = "  +3.43%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'161.22"
$ws.Range("E32").Value = "  +2.74%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.26"
$ws.Range("E33").Value = "  +1.31%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.0744"
$ws.Range("E35").Value = "  +0.89%  "

# Row 36 - LidoDAOToken
$ws.Range("D36").Value = "'3.01"
$ws.Range("E36").Value = "  -1.42%  "

# Row 37 - Celestia
$ws.Range("D37").Value = "'17.04"
$ws.Range("E37").Value = "  +2.51%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +1.93%  "

# Row 39 - WEMIXToken
$ws.Range("E39").Value = "  -0.24%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.62%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +1.43%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "'3.92"
$ws.Range("E42").Value = "  -1.78%  "

# Row 43 & 44 - swap Maker and EnergySwap (with updated values)
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.67"
$ws.Range("E43").Value = "  -3.56%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.002.74"
$ws.Range("E44").Value = "  -3.19%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +1.96%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "'10.31"
$ws.Range("E46").Value = "  +1.96%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +3.04%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'2.90"
$ws.Range("E48").Value = "  -2.86%  "

# Row 49 - MultiversX
$ws.Range("D49").Value = "'53.03"
$ws.Range("E49").Value = "  +3.12%  "

# Row 50 - Stacks
$ws.Range("D50").Value = "'1.51"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  +0.38%  "
